$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update data rows 2-33 (asset table) to reflect refreshed portfolio data
$data = @(
    ,("BTC", 0.006, 45000, "Digital Gold")
    ,("ETH", 6.458, 2600, "Layer 1")
    ,("SOL", 42.31, 40, "Layer 1")
    ,("VET", 186842, 0.0211, "Utility")
    ,("ATOM", 1239.6, 8.5, "Utility")
    ,("LINK", 557.8, 8, "Oracle")
    ,("QNT", 11.26, 83, "Utility")
    ,("MATIC", 2315.7, 1.1, "Layer 2")
    ,("ONT", 1643, 1.1, "Utility")
    ,("AAVE", 49.1, 66.3, "DeFi")
    ,("OCEAN", 5746.13, 0.351, "Storage")
    ,("FET", 10890, 0.375, "KI")
    ,("PYR", 364.1, 4.424, "Gaming")
    ,("HIGH", 1182.2, 2.85, "Metaverse")
    ,("VANRY", 56240, 0.08, "Metaverse")
    ,("AR", 27.14, 31, "Storage")
    ,("AGLD", 228, 0.494, "Gaming")
    ,("API3", 105.44, 9, "Oracle")
    ,("TRB", 22.74, 44, "Oracle")
    ,("SUPER", 2667, 0.68, "Metaverse")
    ,("VTHO", 441238, 0.00255, "Utility")
    ,("POWR", 1191.8, 0.86, "Utility")
    ,("PHB", 2341.8, 1, "KI")
    ,("DESO", 42.56, 46, "DeSo")
    ,("PRIME", 68.66, 8, "Gaming")
    ,("OSMO", 941.6, 1.7, "DeFi")
    ,("TIA", 77.5, 19, "Utility")
    ,("AGIX", 1232, 0.435, "KI")
    ,("INJ", 30.67, 37.6, "DeFi")
    ,("KUJI", 52, 3.99, "Utility")
    ,("MDT", 3884, 0.123, "KI")
    ,("FORT", 1787, 0.307, "KI")
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

# Update sheet view: scroll position and active selection
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("E33").Select()
